# update on 20210624 孤岛风云
# Add Korean (ko_KR) translations for the two story lines in column D.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "지마는 악몽에서 깨어날 수 없었다.`n"
$ws.Range("D3").Value = "그녀는 늘 이런 꿈을 꾼다. 그녀는 언제까지나 계속 이런 꿈을 꿀 것이다.`n"

# Re-run autofit so the newline we just introduced doesn't leave a stray
# explicit row-height behind (matches the source diff, which only touches
# the shared strings table and the two cell references).
$ws.Rows(2).EntireRow.AutoFit()
$ws.Rows(3).EntireRow.AutoFit()
